# Re-run SGNN to annotate dialog acts following clean up work to the
# original transcripts. Updates columns I (DAMSLTag) and J (DialogAct)
# for the rows whose annotations changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new DAMSLTag / DialogAct values
$updates = @(
    @{ Row = 4;  DAMSLTag = "sv"; DialogAct = "Statement-opinion" },
    @{ Row = 14; DAMSLTag = "ba"; DialogAct = "Appreciation" },
    @{ Row = 15; DAMSLTag = "ba"; DialogAct = "Appreciation" },
    @{ Row = 17; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 22; DAMSLTag = "ba"; DialogAct = "Appreciation" },
    @{ Row = 29; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 38; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 66; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 72; DAMSLTag = "aa"; DialogAct = "Agree/Accept" },
    @{ Row = 74; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 83; DAMSLTag = "%";  DialogAct = "Uninterpretable" },
    @{ Row = 84; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 95; DAMSLTag = "sv"; DialogAct = "Statement-opinion" },
    @{ Row = 99; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.DAMSLTag
    $ws.Cells.Item($u.Row, 10).Value = $u.DialogAct
}
